$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 6): waldensian / (0, 210, 255)
$ws.Range("A6").Value = "waldensian"
$ws.Range("B6").Value = "(0, 210, 255)"

# Update the selection to match the saved view state
$ws.Range("E4").Select() | Out-Null
